# Apply "Add data for 2023-08-24" cell updates across the affected sheets.
$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet 1)
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("H2").Value = 70
$ws.Range("D3").Value = 93
$ws.Range("B6").Value = 251
$ws.Range("C6").Value = 321
$ws.Range("E6").Value = 284
$ws.Range("F6").Value = 375
$ws.Range("G6").Value = 327
$ws.Range("H6").Value = 291
$ws.Range("I6").Value = 362
$ws.Range("J6").Value = 271
$ws.Range("B7").Value = 343
$ws.Range("C7").Value = 432
$ws.Range("D7").Value = 450
$ws.Range("E7").Value = 433
$ws.Range("F7").Value = 532
$ws.Range("G7").Value = 477
$ws.Range("H7").Value = 462
$ws.Range("I7").Value = 594
$ws.Range("J7").Value = 506

# Garfield Park (sheet 10)
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("B6").Value = 9
$ws.Range("I6").Value = 23
$ws.Range("B7").Value = 11
$ws.Range("I7").Value = 35

# Grand Crossing (sheet 11)
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C6").Value = 23
$ws.Range("C7").Value = 28

# Humboldt Park (sheet 15)
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J4").Value = 6
$ws.Range("J5").Value = 10

# By Neighborhood (sheet 2)
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F8").Value = 34
$ws.Range("D19").Value = 19
$ws.Range("B32").Value = 11
$ws.Range("I32").Value = 35
$ws.Range("C36").Value = 28
$ws.Range("J41").Value = 10
$ws.Range("B53").Value = 33
$ws.Range("C53").Value = 38
$ws.Range("E53").Value = 52
$ws.Range("G53").Value = 60
$ws.Range("H53").Value = 55
$ws.Range("I53").Value = 93
$ws.Range("G76").Value = 10
$ws.Range("J76").Value = 11
$ws.Range("I77").Value = 32
$ws.Range("J77").Value = 25
$ws.Range("C85").Value = 12
$ws.Range("J96").Value = 7
$ws.Range("B98").Value = 343
$ws.Range("C98").Value = 432
$ws.Range("D98").Value = 450
$ws.Range("E98").Value = 433
$ws.Range("F98").Value = 532
$ws.Range("G98").Value = 477
$ws.Range("H98").Value = 462
$ws.Range("I98").Value = 594
$ws.Range("J98").Value = 506

# Woodlawn (sheet 21)
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value = 4
$ws.Range("J6").Value = 7

# Loop (sheet 22)
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("H2").Value = 6
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 23
$ws.Range("E6").Value = 41
$ws.Range("G6").Value = 41
$ws.Range("H6").Value = 38
$ws.Range("I6").Value = 60
$ws.Range("B7").Value = 33
$ws.Range("C7").Value = 38
$ws.Range("E7").Value = 52
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 93

# Rogers Park (sheet 3)
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("G5").Value = 5
$ws.Range("J5").Value = 9
$ws.Range("G6").Value = 10
$ws.Range("J6").Value = 11

# Roseland (sheet 4)
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I6").Value = 20
$ws.Range("J6").Value = 14
$ws.Range("I7").Value = 32
$ws.Range("J7").Value = 25

# Chatham (sheet 47)
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("D3").Value = 4
$ws.Range("D6").Value = 19

# United Center (sheet 52)
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("C4").Value = 9
$ws.Range("C5").Value = 12

# Austin (sheet 7)
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F5").Value = 23
$ws.Range("F6").Value = 34
